$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46037
$ws.Range("B2").Value = 119.05
$ws.Range("C2").Value = 107.56
$ws.Range("D2").Value = 101.39
$ws.Range("E2").Value = 94.83
$ws.Range("F2").Value = 87.04000000000001
$ws.Range("G2").Value = 86.31
$ws.Range("H2").Value = 98.90000000000001
$ws.Range("I2").Value = 112.55
$ws.Range("J2").Value = 114.7
$ws.Range("K2").Value = 105.48
$ws.Range("L2").Value = 98.59
$ws.Range("M2").Value = 94.89
$ws.Range("N2").Value = 92.2
$ws.Range("O2").Value = 86.67
$ws.Range("P2").Value = 89.83
$ws.Range("Q2").Value = 94.51000000000001
$ws.Range("R2").Value = 105.03
$ws.Range("S2").Value = 113.14
$ws.Range("T2").Value = 116.66
$ws.Range("U2").Value = 120.72
$ws.Range("V2").Value = 133.07
$ws.Range("W2").Value = 138.2
$ws.Range("X2").Value = 128.64
$ws.Range("Y2").Value = 113.8
$ws.Range("Z2").Value = 106.41
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 128.43
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 135.64
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 121.22
$ws.Range("AG2").Value = "2h-16h"
